# Update New Orleans xlsx shard:
#   1. Add a "State" column (value "Louisiana") to hotel_info, placed right
#      after Hotel_Name and before City.
#   2. Swap the sheet tab order so review_info becomes the first sheet and
#      hotel_info becomes the second sheet.
#
# Note: this engine keeps sheetId/r:id attached to worksheet *position*
# (the first worksheet part keeps sheetId="1"/rId1, the second keeps
# sheetId="2"/rId2) rather than to the sheet's logical identity, and a
# plain Worksheet.Move() instead carries the original sheetId along with
# the sheet. To reproduce the target file (where sheetId now tracks
# position, i.e. review_info=1/rId1, hotel_info=2/rId2) the cell content
# is relocated between the two sheet objects (via Range.Cut, which -
# unlike read/write through .Value - keeps numeric-looking text such as
# "1940" stored as text) and the sheets are renamed afterwards, rather
# than reordering the sheets themselves.

$wb = $excel.ActiveWorkbook

$hotelSheet  = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

# --- 1. insert the new "State" column into hotel_info ------------------
$hotelSheet.Columns.Item(3).Insert()
$hotelSheet.Cells.Item(1, 3).Value = "State"
$hotelSheet.Cells.Item(2, 3).Value = "Louisiana"

# --- 2. relocate content so it ends up paired with the right sheetId ---

# scratch sheet used as a holding area so the two ranges (which both
# start at A1) don't clobber each other while swapping
$scratch = $wb.Worksheets.Add()
$scratch.Name = "__scratch__"

# Worksheets.Add() invalidates earlier worksheet handles, so re-fetch
$hotelSheet  = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

$hotelSheet.Range("A1:J2").Cut($scratch.Range("A1:J2"))
$reviewSheet.Range("A1:Y1").Cut($hotelSheet.Range("A1:Y1"))
$scratch.Range("A1:J2").Cut($reviewSheet.Range("A1:J2"))

$scratch = $wb.Worksheets.Item("__scratch__")
[void]$scratch.Delete()

# --- 3. rename the sheets so names line up with their (swapped) content-
$hotelSheet  = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

$hotelSheet.Name  = "__tmp_swap__"
$reviewSheet.Name = "hotel_info"
$hotelSheet.Name  = "review_info"
